$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# Update marking scheme for "Right" answers (row 11, Marking)
$ws.Range("B11").Value = 5

# Update total score (row 12, Total)
$ws.Range("B12").Value = 115

# Update the "correct/total" summary text
$ws.Range("E12").Value = "115/140"
